$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts Category */Sub Category * to F/G)
$ws.Range("E1").EntireColumn.Insert()

# New column takes on the same width as the Portfolio Company * column to its left
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# Header for the new column
$ws.Range("E1").Value = "Pan *"

# Fill in Pan values per-row, matching the Portfolio Company (XYZ -> A11111111, ABC -> B11111111)
$ws.Range("E2").Value = "A11111111"
$ws.Range("E3").Value = "A11111111"
$ws.Range("E4").Value = "B11111111"
$ws.Range("E5").Value = "B11111111"

# Update selection to match the final saved state
$ws.Range("E6").Select()
